$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.594.77"
$ws.Range("E2").Value = '  +0.64%  '
$ws.Range("D3").Value = "'1.923.07"
$ws.Range("E3").Value = '  -0.18%  '
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = "'247.08"
$ws.Range("E5").Value = '  +2.71%  '
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").Value = "'0.4741"
$ws.Range("E7").Value = '  -0.35%  '
$ws.Range("D8").Value = "'0.2904"
$ws.Range("E8").Value = '  +1.62%  '
$ws.Range("D9").Value = "'0.06831"
$ws.Range("E9").Value = '  +3.99%  '
$ws.Range("D10").Value = "'105.53"
$ws.Range("E10").Value = '  -0.67%  '
$ws.Range("D11").Value = "'18.42"
$ws.Range("E11").Value = '  -3.29%  '
$ws.Range("D12").Value = "'1.927.72"
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").Value = "'0.07698"
$ws.Range("E13").Value = '  +1.32%  '
$ws.Range("D14").Value = "'5.353"
$ws.Range("E14").Value = '  +4.48%  '
$ws.Range("D15").Value = "'0.6711"
$ws.Range("E15").Value = '  +2.46%  '
$ws.Range("D16").Value = "'290.03"
$ws.Range("E16").Value = '  -4.11%  '
$ws.Range("D17").Value = "'30.614.09"
$ws.Range("E17").Value = '  +0.65%  '
$ws.Range("D18").Value = "'0.000007621"
$ws.Range("E18").Value = '  +1.85%  '
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").Value = "'12.95"
$ws.Range("E20").Value = '  +0.19%  '
$ws.Range("D21").Value = "'5.549"
$ws.Range("E21").Value = '  +4.93%  '
$ws.Range("D22").Value = "'2.176.24"
$ws.Range("E22").Value = '  +0.75%  '
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").Value = "'6.494"
$ws.Range("E24").Value = '  +3.61%  '
$ws.Range("D25").Value = "'9.506"
$ws.Range("E25").Value = '  +3.32%  '
$ws.Range("D26").Value = "'167.48"
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("D27").Value = "'21.06"
$ws.Range("E27").Value = '  +5.24%  '
$ws.Range("D28").Value = "'2.121"
$ws.Range("E28").Value = '  +5.06%  '
$ws.Range("D29").Value = "'0.1071"
$ws.Range("E29").Value = '  -3.80%  '
$ws.Range("D30").Value = "'1.404"
$ws.Range("E30").Value = '  +3.73%  '
$ws.Range("D31").Value = "'4.181"
$ws.Range("E31").Value = '  +2.34%  '
$ws.Range("D32").Value = "'4.048"
$ws.Range("E32").Value = '  +3.47%  '
$ws.Range("D33").Value = "'0.05026"
$ws.Range("E33").Value = '  +0.78%  '
$ws.Range("D34").Value = "'0.7321"
$ws.Range("E34").Value = '  -0.97%  '
$ws.Range("D35").Value = "'1.145"
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("D36").Value = "'0.02069"
$ws.Range("E36").Value = '  +6.98%  '
$ws.Range("D37").Value = "'1.000"
$ws.Range("E37").Value = '  +0.08%  '
$ws.Range("D38").Value = "'2.725"
$ws.Range("E38").Value = '  -0.77%  '
$ws.Range("D39").Value = "'2.677"
$ws.Range("E39").Value = '  -0.83%  '
$ws.Range("D40").Value = "'112.00"
$ws.Range("E40").Value = '  +4.85%  '
$ws.Range("D41").Value = "'2.041"
$ws.Range("E41").Value = '  -0.76%  '
$ws.Range("D42").Value = "'0.8731"
$ws.Range("E42").Value = '  -0.47%  '
$ws.Range("D43").Value = "'0.4405"
$ws.Range("E43").Value = '  +6.67%  '
$ws.Range("D44").Value = "'5.922"
$ws.Range("E44").Value = '  +2.07%  '
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("D46").Value = "'67.86"
$ws.Range("E46").Value = '  -3.00%  '
$ws.Range("D47").Value = "'7.291"
$ws.Range("E47").Value = '  +0.77%  '
$ws.Range("D48").Value = "'9.388"
$ws.Range("E48").Value = '  +0.99%  '
$ws.Range("D49").Value = "'48.50"
$ws.Range("E49").Value = '  +15.81%  '
$ws.Range("D50").Value = "'0.1244"
$ws.Range("E50").Value = '  +3.89%  '
$ws.Range("D51").Value = "'35.03"
$ws.Range("E51").Value = '  +0.60%  '
